# Update odds for the Wales - Cymru Premier matches (rows 4, 5, 7) in the
# FlashScore weekly export, matching the latest scraped market prices.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G4").Value = 2.4
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 2.87
$ws.Range("K4").Value = 2.05
$ws.Range("L4").Value = 3.5
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 2.92
$ws.Range("Q4").Value = 2.07
$ws.Range("R4").Value = 1.7
$ws.Range("U4").Value = 1.78
$ws.Range("V4").Value = 1.93
$ws.Range("W4").Value = 7.6
$ws.Range("Y4").Value = 9.25
$ws.Range("Z4").Value = 26
$ws.Range("AA4").Value = 20
$ws.Range("AB4").Value = 30
$ws.Range("AC4").Value = 6.5
$ws.Range("AD4").Value = 5.9
$ws.Range("AE4").Value = 13.5
$ws.Range("AF4").Value = 65
$ws.Range("AG4").Value = 8.5
$ws.Range("AH4").Value = 14.5
$ws.Range("AI4").Value = 10.5
$ws.Range("AJ4").Value = 37
$ws.Range("AK4").Value = 26
$ws.Range("AL4").Value = 35
$ws.Range("AM4").Value = 500
$ws.Range("AP4").Value = 20
$ws.Range("AR4").Value = 80
$ws.Range("AU4").Value = 6.7
$ws.Range("AW4").Value = 4.9
$ws.Range("AX4").Value = 16
$ws.Range("AY4").Value = 23
$ws.Range("AZ4").Value = 75
$ws.Range("BA4").Value = 110
$ws.Range("BB4").Value = 300
$ws.Range("G5").Value = 7.4
$ws.Range("J5").Value = 6.9
$ws.Range("K5").Value = 2.32
$ws.Range("L5").Value = 1.87
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.6
$ws.Range("Q5").Value = 1.75
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 1.36
$ws.Range("T5").Value = 2.9
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.72
$ws.Range("AE5").Value = 20
$ws.Range("AG5").Value = 6.6
$ws.Range("AH5").Value = 6.3
$ws.Range("AK5").Value = 11.5
$ws.Range("AL5").Value = 28
$ws.Range("AT5").Value = 2.9
$ws.Range("AV5").Value = 80
$ws.Range("AX5").Value = 6.3
$ws.Range("AY5").Value = 17
$ws.Range("AZ5").Value = 18
$ws.Range("G7").Value = 4.9
$ws.Range("H7").Value = 3.8
$ws.Range("I7").Value = 1.57
$ws.Range("J7").Value = 5.2
$ws.Range("K7").Value = 2.2
$ws.Range("L7").Value = 2.15
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = 1.24
$ws.Range("P7").Value = 3.6
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 2.02
$ws.Range("S7").Value = 1.38
$ws.Range("T7").Value = 2.8
$ws.Range("U7").Value = 1.78
$ws.Range("V7").Value = 1.93
$ws.Range("W7").Value = 15
$ws.Range("Y7").Value = 16
$ws.Range("Z7").Value = 90
$ws.Range("AA7").Value = 50
$ws.Range("AB7").Value = 50
$ws.Range("AC7").Value = 8
$ws.Range("AD7").Value = 7.6
$ws.Range("AF7").Value = 70
$ws.Range("AG7").Value = 7.3
$ws.Range("AH7").Value = 7.7
$ws.Range("AJ7").Value = 11.75
$ws.Range("AL7").Value = 24
$ws.Range("AM7").Value = 500
$ws.Range("AO7").Value = 29
$ws.Range("AP7").Value = 35
$ws.Range("AQ7").Value = 175
$ws.Range("AR7").Value = 200
$ws.Range("AS7").Value = 450
$ws.Range("AT7").Value = 2.8
$ws.Range("AU7").Value = 7.7
$ws.Range("AV7").Value = 75
$ws.Range("AW7").Value = 3.45
$ws.Range("AX7").Value = 7.8
$ws.Range("AY7").Value = 17.5
$ws.Range("AZ7").Value = 25
